# Add a new record (row 8) to the inspection log on Sheet1.
# The sheet currently has data in rows 1-7 (A1:R7); this appends one more
# data row, extending the used range to A1:R8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-looking value to be stored as literal
# text (matching the other rows, which keep "Fecha" values as plain text
# rather than auto-converted date serials).
$ws.Range("A8").Value = "'2024-08-23"
$ws.Range("B8").Value = "II"
$ws.Range("C8").Value = "jfsdjfldsk"
$ws.Range("D8").Value = "BUSTAMANTE MERCHAN RONALD ALEJANDRO"
$ws.Range("E8").Value = "Cumple"
# F8:L8 and O8 are left blank (no value), matching the source row.
$ws.Range("M8").Value = "Cumple"
$ws.Range("N8").Value = "No cumple"
$ws.Range("P8").Value = "Cumple"
$ws.Range("Q8").Value = "PABLO ENRIQUEZ"
$ws.Range("R8").Value = "vdjxlkvjcxlk"
